# Edit script: rename "Requested quantity" headers and add a new
# "PO Forecast" sheet with forecast data (ds, PO_Forecast, yhat_lower, yhat_upper).

$wb = $excel.ActiveWorkbook

# --- 1) Rename column headers on existing sheets ---
$wsWeekly = $wb.Worksheets.Item(1)   # "Weekly Quantity"
$wsWeekly.Range("B1").Value = "Weekly_PO_Qty"

$wsMonthly = $wb.Worksheets.Item(2)  # "Monthly Trend"
$wsMonthly.Range("B1").Value = "Monthly_PO_Qty"

# --- 2) Add the new "PO Forecast" sheet after the last existing sheet ---
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws3 = $wb.Worksheets.Add($null, $lastSheet)
$ws3.Name = "PO Forecast"

# Copy header formatting (bold/centered/bordered) from the Weekly Quantity header row
$wsWeekly.Range("A1:B1").Copy()
$ws3.Range("A1:D1").PasteSpecial(-4122)

# Copy the date-number-format style down column A for the 28 data rows
$wsWeekly.Range("A2").Copy()
$ws3.Range("A2:A29").PasteSpecial(-4122)

# --- 3) Set header labels ---
$ws3.Range("A1").Value = "ds"
$ws3.Range("B1").Value = "PO_Forecast"
$ws3.Range("C1").Value = "yhat_lower"
$ws3.Range("D1").Value = "yhat_upper"

# --- 4) Fill forecast data rows ---
$rowIndex = 1
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 44934.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 4; $ws3.Cells.Item($rowIndex, 3).Value = -0.2591803429840951; $ws3.Cells.Item($rowIndex, 4).Value = 8.538467453613412
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 44941.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 4; $ws3.Cells.Item($rowIndex, 3).Value = -0.1776220099695499; $ws3.Cells.Item($rowIndex, 4).Value = 8.906405469197727
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 44955.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 4; $ws3.Cells.Item($rowIndex, 3).Value = -0.7917602603296198; $ws3.Cells.Item($rowIndex, 4).Value = 8.421411503201538
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 44962.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 4; $ws3.Cells.Item($rowIndex, 3).Value = -0.8747254470783714; $ws3.Cells.Item($rowIndex, 4).Value = 8.272039441644475
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 44969.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 4; $ws3.Cells.Item($rowIndex, 3).Value = -0.794427648418226; $ws3.Cells.Item($rowIndex, 4).Value = 8.594807935505775
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 44976.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 4; $ws3.Cells.Item($rowIndex, 3).Value = -0.7346544759304557; $ws3.Cells.Item($rowIndex, 4).Value = 7.880174521441363
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 44983.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 4; $ws3.Cells.Item($rowIndex, 3).Value = -1.012754690526311; $ws3.Cells.Item($rowIndex, 4).Value = 7.864844323321598
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 44997.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 3; $ws3.Cells.Item($rowIndex, 3).Value = -1.384633520081198; $ws3.Cells.Item($rowIndex, 4).Value = 8.22284152779681
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45011.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 3; $ws3.Cells.Item($rowIndex, 3).Value = -1.347043570018074; $ws3.Cells.Item($rowIndex, 4).Value = 7.819854532606583
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45018.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 3; $ws3.Cells.Item($rowIndex, 3).Value = -1.160530707910143; $ws3.Cells.Item($rowIndex, 4).Value = 8.055276097215479
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45025.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 3; $ws3.Cells.Item($rowIndex, 3).Value = -1.339459287406098; $ws3.Cells.Item($rowIndex, 4).Value = 7.675296438633049
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45032.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 3; $ws3.Cells.Item($rowIndex, 3).Value = -1.350112910714259; $ws3.Cells.Item($rowIndex, 4).Value = 7.66700234092318
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45039.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 3; $ws3.Cells.Item($rowIndex, 3).Value = -1.649760487220139; $ws3.Cells.Item($rowIndex, 4).Value = 7.151312696809547
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45046.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 3; $ws3.Cells.Item($rowIndex, 3).Value = -1.589445206700117; $ws3.Cells.Item($rowIndex, 4).Value = 7.312858963965952
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45053.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 3; $ws3.Cells.Item($rowIndex, 3).Value = -1.864231622278816; $ws3.Cells.Item($rowIndex, 4).Value = 7.133931826474719
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45060.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 3; $ws3.Cells.Item($rowIndex, 3).Value = -1.882377226143845; $ws3.Cells.Item($rowIndex, 4).Value = 6.953962892227002
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45067.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 3; $ws3.Cells.Item($rowIndex, 3).Value = -2.092932240988755; $ws3.Cells.Item($rowIndex, 4).Value = 7.19948393297773
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45074.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 2; $ws3.Cells.Item($rowIndex, 3).Value = -2.030059328613417; $ws3.Cells.Item($rowIndex, 4).Value = 7.213096630477263
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45081.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 2; $ws3.Cells.Item($rowIndex, 3).Value = -2.527939724424794; $ws3.Cells.Item($rowIndex, 4).Value = 6.86089547201718
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45088.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 2; $ws3.Cells.Item($rowIndex, 3).Value = -2.341701239410636; $ws3.Cells.Item($rowIndex, 4).Value = 6.996135709774709
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45095.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 2; $ws3.Cells.Item($rowIndex, 3).Value = -2.417153194384885; $ws3.Cells.Item($rowIndex, 4).Value = 6.949460370735401
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45102.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 2; $ws3.Cells.Item($rowIndex, 3).Value = -2.53654589771281; $ws3.Cells.Item($rowIndex, 4).Value = 6.297020087707792
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45109.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 2; $ws3.Cells.Item($rowIndex, 3).Value = -2.54506785920033; $ws3.Cells.Item($rowIndex, 4).Value = 6.286596146613912
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45116.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 2; $ws3.Cells.Item($rowIndex, 3).Value = -2.48173283100756; $ws3.Cells.Item($rowIndex, 4).Value = 6.452605574515645
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45123.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 2; $ws3.Cells.Item($rowIndex, 3).Value = -2.53232773543358; $ws3.Cells.Item($rowIndex, 4).Value = 6.621038110917127
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45130.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 2; $ws3.Cells.Item($rowIndex, 3).Value = -2.891184783419642; $ws3.Cells.Item($rowIndex, 4).Value = 6.320840137056916
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45137.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 2; $ws3.Cells.Item($rowIndex, 3).Value = -3.002713376351964; $ws3.Cells.Item($rowIndex, 4).Value = 6.316039010541536
$rowIndex++; $ws3.Cells.Item($rowIndex, 1).Value = 45144.99999999999; $ws3.Cells.Item($rowIndex, 2).Value = 2; $ws3.Cells.Item($rowIndex, 3).Value = -2.95447298193994; $ws3.Cells.Item($rowIndex, 4).Value = 6.437158856447184
